# Append a new data row (row 95) to the worksheet, mirroring the existing
# Adafruit IO feed data rows (Timestamp, Feed Key, Value, Latitude,
# Longitude, Elevation), all stored as text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 95

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# Column C holds a numeric-looking reading ("25") that must stay text,
# matching the rest of the sheet (every value in this feed export is
# stored as text). A leading apostrophe is the normal Excel way of
# forcing a numeric-looking entry to be kept as text.
$ws.Cells.Item($row, 3).Value = "'25"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
